# Generate Report for Handback
# Applies the localization-status.xlsx report refresh:
#  - Status text updated from "Ready for handoff" to "Handed back: in sync with en-US"
#    (Overview!E2/F2, zh-cn!C2, de-de!C2)
#  - Latest Handback DateTime refreshed for zh-cn and de-de
#  - Error Detail (stale-handback warning) cleared now that both locales are in sync
#  - A few report columns were widened/narrowed to better fit the refreshed content

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# --- Status text: "Ready for handoff" -> "Handed back: in sync with en-US" ---
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("C2").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("C2").Value = "Handed back: in sync with en-US"

# --- Latest Handback DateTime refreshed ---
$wsZhCn.Range("K2").Value = "2016-09-02 10:58:29"
$wsDeDe.Range("K2").Value = "2016-09-02 10:58:36"

# --- Error Detail cleared (handback is now in sync, no stale-version warning) ---
$wsZhCn.Range("P2").Value = ""
$wsDeDe.Range("P2").Value = ""

# --- Column width adjustments (report column resize) ---
$wsOverview.Columns.Item(5).ColumnWidth = 16.333333333333336
$wsOverview.Columns.Item(6).ColumnWidth = 16.333333333333336

$wsZhCn.Columns.Item(3).ColumnWidth = 29.166666666666664
$wsZhCn.Columns.Item(16).ColumnWidth = 12.833333333333332

$wsDeDe.Columns.Item(3).ColumnWidth = 29.166666666666664
$wsDeDe.Columns.Item(16).ColumnWidth = 12.833333333333332
